# Commit: "Moved discovery process to start in presentation"
#
# The "Discovery Process" slide used to sit later in the deck (right before
# the "Critical Analysis" slides, i.e. 8th out of 10). Move it so it becomes
# the 3rd slide - right after "Technology overview" and before "Metrics" -
# while leaving every other slide's content/order untouched.

$p = $ppt.ActivePresentation

$targetSlide = $null

# Primary lookup: find the slide whose title is "Discovery Process", so the
# script still works even if the deck's current slide order differs from
# what we expect.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    $title = $candidate.Shapes.Item(1).TextFrame.TextRange.Text
    if ($title -eq "Discovery Process") {
        $targetSlide = $candidate
        break
    }
}

# Fallback lookup: the slide's original, stable SlideID in the source deck.
if ($targetSlide -eq $null) {
    $targetSlide = $p.Slides.FindBySlideID(264)
}

# Move it to position 3.
$targetSlide.MoveTo(3)
